# Weekly price-sheet update: a new weekly record is inserted as row 407
# (pushing the existing rows 407-433 down to 408-434). Implemented as an
# explicit "insert row" by shifting the variable columns (D,J,K,L,M,O,P)
# down by one row, working from the bottom up so we never clobber data
# before it has been read, and then populating the brand-new row 407.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 407
$lastRow = 433
$newLastRow = $lastRow + 1

# Make sure the new last row shares the date-column number format with
# the rest of the column (it's a brand-new row, so it has no formatting
# of its own yet).
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat

# Row $newLastRow doesn't exist yet, so the columns that stay constant
# across the whole block (A,B,C,E,F,G,H,I,N,Q,R) need to be copied down
# explicitly too - the shift loop below only moves the columns that vary
# row to row (D,J,K,L,M,O,P).
foreach ($col in @("A","B","C","E","F","G","H","I","N","Q","R")) {
    $ws.Range("$col$newLastRow").Value = $ws.Range("$col$lastRow").Value2
}

# Shift rows down: new row (r+1) gets the values that used to live in row r.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $dest = $r + 1

    $dVal = $ws.Range("D$r").Value2
    $jVal = $ws.Range("J$r").Value2
    $kVal = $ws.Range("K$r").Value2
    $lVal = $ws.Range("L$r").Value2
    $mVal = $ws.Range("M$r").Value2
    $oVal = $ws.Range("O$r").Value2
    $pVal = $ws.Range("P$r").Value2

    $ws.Range("D$dest").Value = $dVal
    $ws.Range("J$dest").Value = $jVal
    $ws.Range("K$dest").Value = $kVal
    $ws.Range("L$dest").Value = $lVal
    $ws.Range("M$dest").Value = $mVal
    $ws.Range("O$dest").Value = $oVal
    $ws.Range("P$dest").Value = $pVal
}

# Populate the brand-new row 407 with the new weekly record.
$ws.Range("D$firstRow").Value = 44931
$ws.Range("J$firstRow").Value = 80
$ws.Range("K$firstRow").Value = 8000
$ws.Range("L$firstRow").Value = 8000
$ws.Range("M$firstRow").Value = 8000
$ws.Range("O$firstRow").Value = "Provincia de Cautín"
$ws.Range("P$firstRow").Value = 667
